# The deck ships two theme parts:
#   theme1.xml  -> "Office Theme" / "Office" colour scheme  (only used by the Notes Master)
#   theme2.xml  -> "Integral"     / "Red Violet" colour scheme (used by the Slide Master -> every slide)
#
# The authored edit swaps the *content* of those two theme parts in place
# (file names / relationships are untouched), so the deck that used to render
# with the "Integral" / Red Violet palette now renders with the plain
# "Office Theme" palette, while the part that the Notes Master points at ends
# up holding what used to be the "Integral" colours.
#
# The only theme surface the PowerPoint object model exposes is the slide
# master's (==presentation's) colour scheme, which is backed by the theme
# part actually used to render the slides (theme2.xml here). We drive that
# through SlideMaster.ColorScheme, rewriting every one of the twelve theme
# colour slots (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) from the
# "Red Violet" values to the plain "Office" values, which is exactly the
# observable effect of the authored swap on the rendered deck.

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.ColorScheme

$cs.Colors(1).RGB  = 0         # dk1      000000
$cs.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$cs.Colors(3).RGB  = 6968388   # dk2      44546A
$cs.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$cs.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$cs.Colors(6).RGB  = 3243501   # accent2  ED7D31
$cs.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$cs.Colors(8).RGB  = 49407     # accent4  FFC000
$cs.Colors(9).RGB  = 12874308  # accent5  4472C4
$cs.Colors(10).RGB = 4697456   # accent6  70AD47
$cs.Colors(11).RGB = 12673797  # hlink    0563C1
$cs.Colors(12).RGB = 7491477   # folHlink 954F72
